# Goods.xlsx update:
#  - Row 4 (was id=3, 鼠标/罗技) becomes id=4, 水杯/青鸟, price 69 -> 25
#  - Row 5 (was id=4, 水杯/青鸟) becomes id=5, 抽纸/清风, price 25 -> 5, number 50 -> 200
#  - Row 2 number 99 -> 96
#  - Row 3 price 299 -> 199, number 100 -> 97

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: number 99 -> 96
$ws.Range("E2").Value = 96

# Row 3: price 299 -> 199, number 100 -> 97
$ws.Range("D3").Value = 199
$ws.Range("E3").Value = 97

# Row 4: id 3 -> 4, goodsname 鼠标 -> 水杯, producer 罗技 -> 青鸟, price 69 -> 25 (number stays 50)
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "水杯"
$ws.Range("C4").Value = "青鸟"
$ws.Range("D4").Value = 25

# Row 5: id 4 -> 5, goodsname 水杯 -> 抽纸, producer 青鸟 -> 清风, price 25 -> 5, number 50 -> 200
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "抽纸"
$ws.Range("C5").Value = "清风"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 200
